$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Mon_Dec__4_15_55_04_2023", "f", 30),
    @("Mon_Dec__4_16_06_35_2023", "f", 30),
    @("Mon_Dec__4_16_09_58_2023", "f", 30),
    @("Mon_Dec__4_16_12_03_2023", "f", 30),
    @("Mon_Dec__4_16_16_40_2023", "f", 30),
    @("Mon_Dec__4_16_16_43_2023", "f", 30),
    @("Mon_Dec__4_16_17_36_2023", "f", 30)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
